# Auto-generated edit script applying the diff to cryptos sheet (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.442.72'
$ws.Range("E2").Value = '  +2.56%  '
$ws.Range("D3").Value = '2.366.76'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +3.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.16'
$ws.Range("E6").Value = '  +2.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.56'
$ws.Range("E7").Value = '  +7.95%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.551'
$ws.Range("E9").Value = '  +20.03%  '
$ws.Range("E10").Value = '  +7.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.13'
$ws.Range("E11").Value = '  +13.84%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '2.715.34'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.88'
$ws.Range("E14").Value = '  +8.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.79'
$ws.Range("E15").Value = '  +8.80%  '
$ws.Range("E16").Value = '  +7.47%  '
$ws.Range("D17").Value = '2.367.68'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '44.456.27'
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("E19").Value = '  +4.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.51'
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("E21").Value = '  +3.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.62'
$ws.Range("E22").Value = '  +2.52%  '
$ws.Range("E23").Value = '  -3.43%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.40'
$ws.Range("E26").Value = '  +4.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.25'
$ws.Range("E27").Value = '  +1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.48'
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("E29").Value = '  +4.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.85'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  +2.45%  '
$ws.Range("E32").Value = '  +5.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0740'
$ws.Range("E33").Value = '  +7.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.20'
$ws.Range("E34").Value = '  +4.14%  '
$ws.Range("E35").Value = '  +3.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.90'
$ws.Range("E36").Value = '  +7.24%  '
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.48'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0272'
$ws.Range("E39").Value = '  +6.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.99'
$ws.Range("E40").Value = '  +10.31%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  +3.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0987'
$ws.Range("E44").Value = '  +4.32%  '
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("E46").Value = '  +12.88%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.48'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.80'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +3.95%  '
$ws.Range("D50").Value = '1.444.36'
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").Value = '2.590.39'
$ws.Range("E51").Value = '  +0.05%  '
